# Commit: "changing document, table attributes to lowerCamelCase"
# The ObjTables header markers in row 1 and row 2 encode attribute names
# (ObjTablesVersion / Type / Id) that are being renamed to lowerCamelCase
# (objTablesVersion / type / id).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws.Range("A2").Value = "!!ObjTables type='Data' id='TransposedNode'"
